# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6283   # 合肥·首届AT次元时代动漫游戏嘉年华: 6284 -> 6283
$ws1.Range("F3").Value = 568    # 合肥·Holic动漫游戏展: 566 -> 568
$ws1.Range("F7").Value = 344    # 合肥·W·A第五人格同人only2.0: 343 -> 344
$ws1.Range("F8").Value = 1420   # 合肥·第九届环形宇宙动漫游戏嘉年华: 1401 -> 1420

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6283    # 合肥·首届AT次元时代动漫游戏嘉年华: 6284 -> 6283
$ws4.Range("F3").Value = 568     # 合肥·Holic动漫游戏展: 566 -> 568
$ws4.Range("F7").Value = 344     # 合肥·W·A第五人格同人only2.0: 343 -> 344
$ws4.Range("F12").Value = 1420   # 合肥·第九届环形宇宙动漫游戏嘉年华: 1401 -> 1420
